# Fruta / hortaliza, semanal
#
# A new weekly price report row for "Feria Lagunitas de Puerto Montt" / Mango
# is inserted at row 44, pushing the existing rows 44-86 down to 45-87.
#
# This mirrors what Excel does when you right-click a row header and choose
# "Insert" then type the new observation into the now-empty row: formatting
# (the date-style on column D) is inherited from the row above, and the
# worksheet's used-range dimension grows by one row automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44; rows 44:86 shift down to 45:87.
$ws.Rows("44:44").Insert()

# Populate the new row 44 with the latest weekly observation.
$ws.Cells.Item(44, 1).Value  = 4
$ws.Cells.Item(44, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(44, 3).Value  = "Los Lagos"
$ws.Cells.Item(44, 4).Value  = 44484
$ws.Cells.Item(44, 5).Value  = 10
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100108
$ws.Cells.Item(44, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(44, 9).Value  = 100108002
$ws.Cells.Item(44, 10).Value = "Mango"
$ws.Cells.Item(44, 11).Value = "Sin especificar"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 200
$ws.Cells.Item(44, 14).Value = 7500
$ws.Cells.Item(44, 15).Value = 8000
$ws.Cells.Item(44, 16).Value = 7750
$ws.Cells.Item(44, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(44, 18).Value = "Perú"
$ws.Cells.Item(44, 19).Value = 1938
$ws.Cells.Item(44, 20).Value = 4
